# Updated cryptos list (matches GitHub Actions scrape commit).
# Price (D) / Volume(1h) (E) columns are stored as plain text in this sheet
# (no numeric formatting), so numeric-looking D values are written with a
# leading apostrophe to force text, then the style is reset back to
# "Normal" so no stray number-format style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.645.81"
$ws.Range("E2").Value = "  +2.34%  "

$ws.Range("D3").Value = "2.297.48"
$ws.Range("E3").Value = "  +1.45%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'308.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.41%  "

$ws.Range("D6").Value = "'97.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.84%  "

$ws.Range("D7").Value = "'0.533"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.80%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.498"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.94%  "

$ws.Range("D10").Value = "'36.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +13.24%  "

$ws.Range("E12").Value = "  -1.36%  "

$ws.Range("D13").Value = "'6.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.59%  "

$ws.Range("D14").Value = "2.652.84"
$ws.Range("E14").Value = "  +1.48%  "

$ws.Range("E15").Value = "  +3.20%  "

$ws.Range("D16").Value = "2.301.09"
$ws.Range("E16").Value = "  +2.59%  "

$ws.Range("E17").Value = "  +5.61%  "

$ws.Range("D18").Value = "42.545.65"
$ws.Range("E18").Value = "  +2.30%  "

$ws.Range("D19").Value = "'12.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.32%  "

$ws.Range("D20").Value = "0.0₃0921"
$ws.Range("E20").Value = "  +1.95%  "

$ws.Range("D21").Value = "'6.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.17%  "

$ws.Range("D22").Value = "'67.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.97%  "

$ws.Range("D23").Value = "'243.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "

$ws.Range("E24").Value = "  +0.61%  "

$ws.Range("E25").Value = "  +2.38%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "'24.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.57%  "

$ws.Range("D28").Value = "'36.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.25%  "

$ws.Range("D29").Value = "'9.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("E30").Value = "  +2.39%  "

$ws.Range("D31").Value = "'161.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.49%  "

$ws.Range("E32").Value = "  +2.52%  "

$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("E34").Value = "  +1.38%  "

$ws.Range("D35").Value = "'3.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.86%  "

$ws.Range("D36").Value = "'17.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.74%  "

$ws.Range("E37").Value = "  +3.34%  "

$ws.Range("E38").Value = "  +5.41%  "

$ws.Range("D39").Value = "'2.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.57%  "

$ws.Range("E40").Value = "  -0.20%  "

$ws.Range("D41").Value = "'4.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.03%  "

$ws.Range("D42").Value = "'2.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +19.76%  "

$ws.Range("D43").Value = "2.010.35"
$ws.Range("E43").Value = "  -1.81%  "

$ws.Range("D44").Value = "'19.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.22%  "

$ws.Range("E45").Value = "  +3.28%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'10.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.70%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.31%  "

$ws.Range("D48").Value = "'54.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.54%  "

$ws.Range("E49").Value = "  +1.26%  "

$ws.Range("D50").Value = "'72.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("E51").Value = "  -0.66%  "
